# KSA_Cities.xlsx - "Add files via upload" commit
#
# 1) The Arabic name for "Wadi Al Dawasir" (row 91, column C) had a
#    misspelling (وادى الداوسر) - correct it to وادى الدواسر.
# 2) Append three new city rows at the bottom of the table (158-160):
#       Fayfa / فيفاء               (Jazan region / south of Kingdom)
#       Ahad Al Masarhah / أحد المسارحة  (Jazan region / south of Kingdom)
#       Al Mahani / المحانى          (Makkah region / west of Kingdom)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) fix the Arabic spelling -----------------------------------------
$ws.Range("C91").Value = "وادى الدواسر"

# --- 2) append the new rows, reusing row 157's formatting ----------------
$ws.Range("A157:G157").Copy()
$ws.Range("A158:G160").PasteSpecial(-4122)

$ws.Range("A158").Value = "Fayfa"
$ws.Range("B158").Value = "Fayfa"
$ws.Range("C158").Value = "فيفاء"
$ws.Range("D158").Value = 17.246535999999999
$ws.Range("E158").Value = 43.107962000000001
$ws.Range("F158").Value = "منطقة جازان"
$ws.Range("G158").Value = "جنوب المملكة"

$ws.Range("A159").Value = "Ahad Al Masarhah"
$ws.Range("B159").Value = "Ahad Al Masarhah"
$ws.Range("C159").Value = "أحد المسارحة"
$ws.Range("D159").Value = 16.711358000000001
$ws.Range("E159").Value = 42.956038999999997
$ws.Range("F159").Value = "منطقة جازان"
$ws.Range("G159").Value = "جنوب المملكة"

$ws.Range("A160").Value = "Al Mahani"
$ws.Range("B160").Value = "Al Mahani"
$ws.Range("C160").Value = "المحانى"
$ws.Range("D160").Value = 22.497288999999999
$ws.Range("E160").Value = 40.442901999999997
$ws.Range("F160").Value = "منطقة مكة المكرمة"
$ws.Range("G160").Value = "غرب المملكة"

# --- keep the sheet's "whole table" selection in sync with the new size --
$ws.Range("A1:G160").Select() | Out-Null
